$d = $word.ActiveDocument

# --- Merge split runs in header paragraphs into single runs (text-only, no formatting) ---
$d.Content.Find.Execute("Questions: Trigonometric identities (radians)", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Trigonometric identities (radians)", 2) | Out-Null
$d.Content.Find.Execute("Dzhemma Ruseva", $false, $false, $false, $false, $false, $true, 1, $false, "Dzhemma Ruseva", 2) | Out-Null
$d.Content.Find.Execute("A selection of questions on trigonometric identities, where angles are measured in radians.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions on trigonometric identities, where angles are measured in radians.", 2) | Out-Null

# --- Fix m:dPr child element order (sepChr before endChr) inside each OMath zone that has delimiters ---
$d.OMaths.Item(1).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:t>2</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>6</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>4</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d></m:oMath>')
$d.OMaths.Item(2).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:t>10</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>7</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>14</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d></m:oMath>')
$d.OMaths.Item(3).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:t>5</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>6</m:t></m:r></m:num><m:den><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>csc</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:den></m:f></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>15</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r></m:num><m:den><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sec</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:den></m:f></m:e></m:d></m:oMath>')
$d.OMaths.Item(4).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>4</m:t></m:r><m:r><m:t>s</m:t></m:r><m:r><m:t>i</m:t></m:r><m:sSup><m:e><m:r><m:t>n</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(5).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(6).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(7).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(8).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(9).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(10).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>csc</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sec</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(11).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(12).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>2</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:num><m:den><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>1</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:e></m:d></m:den></m:f></m:oMath>')
$d.OMaths.Item(13).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:num><m:den><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>7</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:den></m:f></m:oMath>')
$d.OMaths.Item(14).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:num><m:den><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>θ</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>θ</m:t></m:r></m:e></m:d></m:den></m:f></m:oMath>')
$d.OMaths.Item(15).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>7</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>6</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(16).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(17).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>4</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(18).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>18</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(21).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>13</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>18</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(22).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(23).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(24).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(25).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(26).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')
$d.OMaths.Item(27).Range.InsertXML('<m:oMath xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:r><m:t>π</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>/</m:t></m:r><m:r><m:t>12</m:t></m:r></m:e></m:d></m:oMath>')

Write-Output "done; changed=25"
